$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 244, pushing the existing rows 244:345 down to 245:346.
$ws.Rows("244:244").Insert()

# Populate the newly inserted row 244 with this week's new data point.
$ws.Cells.Item(244, 1).Value = 3
$ws.Cells.Item(244, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(244, 3).Value = 'Coquimbo'
$ws.Cells.Item(244, 4).Value = 44704
$ws.Cells.Item(244, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(244, 5).Value = 5
$ws.Cells.Item(244, 6).Value = 100112009
$ws.Cells.Item(244, 7).Value = 'Acelga'
$ws.Cells.Item(244, 8).Value = 'Sin especificar'
$ws.Cells.Item(244, 9).Value = 'Primera'
$ws.Cells.Item(244, 10).Value = 270
$ws.Cells.Item(244, 11).Value = 3000
$ws.Cells.Item(244, 12).Value = 3300
$ws.Cells.Item(244, 13).Value = 3167
$ws.Cells.Item(244, 14).Value = '$/docena de atados (6 kilos)'
$ws.Cells.Item(244, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(244, 16).Value = 528
$ws.Cells.Item(244, 17).Value = 6
$ws.Cells.Item(244, 18).Value = 'Hortaliza'
